$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.803.84"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3
$ws.Range("D3").Value = "2.549.22"
$ws.Range("E3").Value = "  +5.77%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.43"
$ws.Range("E5").Value = "  +2.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.62"
$ws.Range("E6").Value = "  +7.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.80%  "

# Row 9
$ws.Range("D9").Value = "2.548.02"
$ws.Range("E9").Value = "  +5.79%  "

# Row 10
$ws.Range("E10").Value = "  +2.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("E11").Value = "  +0.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("E12").Value = "  +1.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("E13").Value = "  +3.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.17"
$ws.Range("E14").Value = "  +9.64%  "

# Row 15
$ws.Range("D15").Value = "3.005.43"
$ws.Range("E15").Value = "  +5.66%  "

# Row 16
$ws.Range("D16").Value = "63.605.60"
$ws.Range("E16").Value = "  +2.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000144"
$ws.Range("E17").Value = "  +3.19%  "

# Row 18
$ws.Range("D18").Value = "2.550.33"
$ws.Range("E18").Value = "  +5.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.59"
$ws.Range("E19").Value = "  +4.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.61"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +3.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.31"
$ws.Range("E24").Value = "  +2.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -0.69%  "

# Row 26
$ws.Range("E26").Value = "  +3.20%  "

# Row 27
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.36"
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("E29").Value = "  +3.73%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0837"
$ws.Range("E30").Value = "  +7.33%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.92"
$ws.Range("E31").Value = "  +8.75%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.89"
$ws.Range("E32").Value = "  +4.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.60"
$ws.Range("E33").Value = "  +3.43%  "

# Row 34
$ws.Range("E34").Value = "  +14.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "419.18"
$ws.Range("E35").Value = "  +11.57%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.409"
$ws.Range("E36").Value = "  +3.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.19"
$ws.Range("E37").Value = "  +3.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  -1.56%  "

# Row 40
$ws.Range("E40").Value = "  +6.47%  "

# Row 41
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.62"
$ws.Range("E42").Value = "  +3.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "153.50"
$ws.Range("E43").Value = "  +6.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.80"
$ws.Range("E44").Value = "  +3.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.07"
$ws.Range("E45").Value = "  +1.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.613"
$ws.Range("E46").Value = "  +4.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0969"
$ws.Range("E48").Value = "  +0.98%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.94"
$ws.Range("E49").Value = "  +5.53%  "

# Row 50
$ws.Range("E50").Value = "  +5.67%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0234"
$ws.Range("E51").Value = "  +8.65%  "
